$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map old "statut" emoji values to new ones, and update the matching label.
$map = @{
    "🟩" = @{ icon = "📗"; label = "vert" }
    "🟧" = @{ icon = "📙"; label = "orange" }
    "🟥" = @{ icon = "📕"; label = "rouge" }
    "⬛" = @{ icon = "📘"; label = "bleu" }
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    $cellA = $ws.Cells.Item($r, 1)
    $val = $cellA.Value()
    if ($map.ContainsKey($val)) {
        $entry = $map[$val]
        $cellA.Value = $entry.icon
        $ws.Cells.Item($r, 2).Value = $entry.label
    }
}
